$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# India (row 7): update totals
$ws.Range("D7").Value = 147195
$ws.Range("E7").Value = 142587

# Afganistan (row 43): update totals
$ws.Range("B43").Value = 23529
$ws.Range("C43").Value = 639
$ws.Range("D43").Value = 3927
$ws.Range("E43").Value = 19156
$ws.Range("G43").Value = 20
$ws.Range("H43").Value = 446

# Israel now ranks above Panama: row 47 becomes Israel (new data), row 48 becomes Panama (old Panama data)
$ws.Range("A47").Value = "Israel"
$ws.Range("B47").Value = 18701
$ws.Range("C47").Value = 132
$ws.Range("D47").Value = 15288
$ws.Range("E47").Value = 3113
$ws.Range("H47").Value = 300

$ws.Range("A48").Value = "Panama"
$ws.Range("B48").Value = 18586
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 10977
$ws.Range("E48").Value = 7191
$ws.Range("H48").Value = 418

# Australia (row 71): update totals
$ws.Range("B71").Value = 7289
$ws.Range("C71").Value = 4
$ws.Range("D71").Value = 6781
$ws.Range("E71").Value = 406

# Uzbekistan now ranks above Senegal: row 75 becomes Uzbekistan (new data), row 76 becomes Senegal (old Senegal data)
$ws.Range("A75").Value = "Uzbekistan"
$ws.Range("B75").Value = 4819
$ws.Range("C75").Value = 78
$ws.Range("D75").Value = 3637
$ws.Range("E75").Value = 1163
$ws.Range("H75").Value = 19

$ws.Range("A76").Value = "Senegal"
$ws.Range("B76").Value = 4759
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 2994
$ws.Range("E76").Value = 1710
$ws.Range("H76").Value = 55

# Tailandia (row 89): update totals
$ws.Range("B89").Value = 3129
$ws.Range("C89").Value = 4
$ws.Range("E89").Value = 84
